# إضافة حدث جديد في Card16
# Updates row 14 (fills the previously-blank B:K and M cells with the
# literal text "nan", matching the sheet's existing convention for empty
# data cells) and appends a new row 15 with a new service event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# --- Row 14: fill blank cells B14:K14 and M14 with "nan" -------------------
$ws.Range("B14").Value = "nan"
$ws.Range("C14").Value = "nan"
$ws.Range("D14").Value = "nan"
$ws.Range("E14").Value = "nan"
$ws.Range("F14").Value = "nan"
$ws.Range("G14").Value = "nan"
$ws.Range("H14").Value = "nan"
$ws.Range("I14").Value = "nan"
$ws.Range("J14").Value = "nan"
$ws.Range("K14").Value = "nan"
$ws.Range("M14").Value = "nan"

# --- Row 15: new service event ---------------------------------------------
# A15 mirrors A1:A14 ("16" stored as *text*, not a number) - force text via
# a temporary Text number format, then drop the format back to the default
# style so the cell itself stays unstyled, same as its neighbours.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "16"
$ws.Range("A15").Style = "Normal"

# B15:K15 and M15 stay blank (same as the template used for every other
# row), but still exist as real - empty - cells rather than being left out
# of the row entirely: write a placeholder then clear it back out so the
# cell is materialised, then drop the leftover style.
$blankCols = @("B","C","D","E","F","G","H","I","J","K","M")
foreach ($col in $blankCols) {
    $cell = $ws.Range($col + "15")
    $cell.NumberFormat = "@"
    $cell.Value = "x"
    $cell.Value = ""
    $cell.Style = "Normal"
}

$ws.Range("L15").Value = "25\12\2024"
$ws.Range("N15").Value = "تم سن الفلاتس ومعيارها"
$ws.Range("O15").Value = "الخبير"
